$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 15
$ws_ALC.Range("H15").Value = 1058.7667
$ws_ALC.Range("I15").Value = 1058.7667
$ws_ALC.Range("K15").Value = 3176.300099999999
$ws_ALC.Range("M15").Value = -3007.300099999999

# ALC row 87
$ws_ALC.Range("H87").Value = 0
$ws_ALC.Range("J87").Value = 0
$ws_ALC.Range("N87").Value = 0
$ws_ALC.Range("L87").ClearContents()

# ALC row 90
$ws_ALC.Range("H90").Value = 0
$ws_ALC.Range("J90").Value = 0
$ws_ALC.Range("N90").Value = 0
$ws_ALC.Range("L90").ClearContents()

# ALC row 138
$ws_ALC.Range("H138").Value = 3441.8215
$ws_ALC.Range("J138").Value = 3491
$ws_ALC.Range("L138").Value = 10473
$ws_ALC.Range("N138").Value = -20753

# ARM row 4
$ws_ARM.Range("H4").Value = 1221.5
$ws_ARM.Range("I4").Value = 888
$ws_ARM.Range("K4").Value = 888
$ws_ARM.Range("M4").Value = -772

# ARM row 44
$ws_ARM.Range("H44").Value = 51990
$ws_ARM.Range("J44").Value = 51990
$ws_ARM.Range("L44").Value = 51990
$ws_ARM.Range("N44").Value = -52966

# ARM row 80
$ws_ARM.Range("H80").Value = 19499.5
$ws_ARM.Range("I80").Value = 19499.5
$ws_ARM.Range("K80").Value = 19499.5
$ws_ARM.Range("M80").Value = -18501.5

# ARM row 83
$ws_ARM.Range("H83").Value = 19499.5
$ws_ARM.Range("I83").Value = 19499.5
$ws_ARM.Range("K83").Value = 58498.5
$ws_ARM.Range("M83").Value = -53506.5

# BSM row 82
$ws_BSM.Range("H82").Value = 22639.7
$ws_BSM.Range("I82").Value = 17377.555
$ws_BSM.Range("J82").Value = 69999
$ws_BSM.Range("K82").Value = 17377.555
$ws_BSM.Range("L82").Value = 69999
$ws_BSM.Range("M82").Value = -16994.555
$ws_BSM.Range("N82").Value = -70765

# BSM row 85
$ws_BSM.Range("H85").Value = 22639.7
$ws_BSM.Range("I85").Value = 17377.555
$ws_BSM.Range("J85").Value = 69999
$ws_BSM.Range("K85").Value = 17377.555
$ws_BSM.Range("L85").Value = 69999
$ws_BSM.Range("M85").Value = -16051.555
$ws_BSM.Range("N85").Value = -72651

# CRP row 2
$ws_CRP.Range("H2").Value = 618.7273
$ws_CRP.Range("I2").Value = 1161.6
$ws_CRP.Range("J2").Value = 166.33333
$ws_CRP.Range("K2").Value = 1161.6
$ws_CRP.Range("L2").Value = 166.33333
$ws_CRP.Range("M2").Value = -1048.6
$ws_CRP.Range("N2").Value = -392.33333

# CRP row 22
$ws_CRP.Range("H22").Value = 0
$ws_CRP.Range("I22").Value = 0
$ws_CRP.Range("J22").Value = 0
$ws_CRP.Range("K22").Value = 0
$ws_CRP.Range("N22").Value = 0
$ws_CRP.Range("L22").ClearContents()
$ws_CRP.Range("M22").ClearContents()

# CRP row 41
$ws_CRP.Range("H41").Value = 29500
$ws_CRP.Range("J41").Value = 0
$ws_CRP.Range("L41").Value = 0
$ws_CRP.Range("N41").ClearContents()

# CRP row 60
$ws_CRP.Range("H60").Value = 17894.736
$ws_CRP.Range("I60").Value = 7928.5713
$ws_CRP.Range("J60").Value = 45800
$ws_CRP.Range("K60").Value = 7928.5713
$ws_CRP.Range("L60").Value = 45800
$ws_CRP.Range("M60").Value = -7417.5713
$ws_CRP.Range("N60").Value = -46822

# CRP row 68
$ws_CRP.Range("H68").Value = 71428.57000000001
$ws_CRP.Range("I68").Value = 50000
$ws_CRP.Range("K68").Value = 50000
$ws_CRP.Range("M68").Value = -49251

# CRP row 71
$ws_CRP.Range("H71").Value = 71428.57000000001
$ws_CRP.Range("I71").Value = 50000
$ws_CRP.Range("K71").Value = 150000
$ws_CRP.Range("M71").Value = -146256

# CRP row 74
$ws_CRP.Range("H74").Value = 0
$ws_CRP.Range("J74").Value = 0
$ws_CRP.Range("N74").Value = 0
$ws_CRP.Range("L74").ClearContents()

# CRP row 77
$ws_CRP.Range("H77").Value = 0
$ws_CRP.Range("J77").Value = 0
$ws_CRP.Range("N77").Value = 0
$ws_CRP.Range("L77").ClearContents()

# CRP row 141
$ws_CRP.Range("H141").Value = 237838.67
$ws_CRP.Range("I141").Value = 42500
$ws_CRP.Range("J141").Value = 335508
$ws_CRP.Range("K141").Value = 42500
$ws_CRP.Range("L141").Value = 335508
$ws_CRP.Range("M141").Value = -37320
$ws_CRP.Range("N141").Value = -345868

# CUL row 4
$ws_CUL.Range("H4").Value = 4024202.8
$ws_CUL.Range("I4").Value = 2168788.8
$ws_CUL.Range("K4").Value = 6506366.399999999
$ws_CUL.Range("M4").Value = -6506254.399999999

# CUL row 5
$ws_CUL.Range("H5").Value = 85560.664
$ws_CUL.Range("I5").Value = 2429.3635
$ws_CUL.Range("J5").Value = 1000005
$ws_CUL.Range("K5").Value = 7288.0905
$ws_CUL.Range("L5").Value = 3000015
$ws_CUL.Range("M5").Value = -7176.0905
$ws_CUL.Range("N5").Value = -3000239

# CUL row 29
$ws_CUL.Range("H29").Value = 85
$ws_CUL.Range("J29").Value = 100
$ws_CUL.Range("L29").Value = 300
$ws_CUL.Range("N29").Value = -854

# CUL row 135
$ws_CUL.Range("H135").Value = 85560.664
$ws_CUL.Range("I135").Value = 2429.3635
$ws_CUL.Range("J135").Value = 1000005
$ws_CUL.Range("K135").Value = 21864.2715
$ws_CUL.Range("L135").Value = 9000045
$ws_CUL.Range("M135").Value = -19329.2715
$ws_CUL.Range("N135").Value = -9005115

# CUL row 139
$ws_CUL.Range("H139").Value = 2567.9167
$ws_CUL.Range("I139").Value = 1153.7778
$ws_CUL.Range("K139").Value = 3461.3334
$ws_CUL.Range("M139").Value = 1678.6666

# CUL row 140
$ws_CUL.Range("H140").Value = 527684.6
$ws_CUL.Range("I140").Value = 527684.6
$ws_CUL.Range("K140").Value = 1583053.8
$ws_CUL.Range("M140").Value = -1577873.8

# GSM row 43
$ws_GSM.Range("H43").Value = 2346.9092
$ws_GSM.Range("I43").Value = 2346.9092
$ws_GSM.Range("J43").Value = 0
$ws_GSM.Range("K43").Value = 2346.9092
$ws_GSM.Range("L43").Value = 0
$ws_GSM.Range("N43").Value = -2195.9092
$ws_GSM.Range("M43").ClearContents()

# GSM row 136
$ws_GSM.Range("H136").Value = 52928.438
$ws_GSM.Range("J136").Value = 52928.438
$ws_GSM.Range("L136").Value = 158785.314
$ws_GSM.Range("N136").Value = -163885.314

# LTW row 21
$ws_LTW.Range("H21").Value = 9999
$ws_LTW.Range("J21").Value = 9999
$ws_LTW.Range("L21").Value = 9999
$ws_LTW.Range("N21").Value = -10347

# LTW row 132
$ws_LTW.Range("H132").Value = 3099.75
$ws_LTW.Range("I132").Value = 5199
$ws_LTW.Range("K132").Value = 15597
$ws_LTW.Range("M132").Value = -13067

# WVR row 54
$ws_WVR.Range("H54").Value = 53055.332
$ws_WVR.Range("J54").Value = 58437.25
$ws_WVR.Range("L54").Value = 58437.25
$ws_WVR.Range("N54").Value = -59477.25

# WVR row 74
$ws_WVR.Range("H74").Value = 22164.334
$ws_WVR.Range("I74").Value = 20995.5
$ws_WVR.Range("K74").Value = 20995.5
$ws_WVR.Range("M74").Value = -20059.5

# WVR row 77
$ws_WVR.Range("H77").Value = 22164.334
$ws_WVR.Range("I77").Value = 20995.5
$ws_WVR.Range("K77").Value = 62986.5
$ws_WVR.Range("M77").Value = -58306.5

# WVR row 107
$ws_WVR.Range("H107").Value = 1187.6897
$ws_WVR.Range("I107").Value = 1231.1666
$ws_WVR.Range("K107").Value = 3693.4998
$ws_WVR.Range("M107").Value = -1773.4998
